$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextValue "D2" '67.188.19'
Set-TextValue "E2" '  +0.07%  '
Set-TextValue "D3" '2.482.63'
Set-TextValue "E3" '  +0.56%  '
Set-TextValue "D4" '0.999'
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '585.09'
Set-TextValue "E5" '  +0.29%  '
Set-TextValue "D6" '172.94'
Set-TextValue "E6" '  +3.33%  '
Set-TextValue "E7" '  -0.02%  '
Set-TextValue "E8" '  +0.15%  '
Set-TextValue "D9" '2.482.29'
Set-TextValue "E9" '  +0.51%  '
Set-TextValue "E10" '  +3.05%  '
Set-TextValue "E11" '  +1.09%  '
Set-TextValue "E12" '  +0.35%  '
Set-TextValue "E13" '  +0.10%  '
Set-TextValue "E15" '  +0.39%  '
Set-TextValue "D16" '67.010.42'
Set-TextValue "E16" '  +0.35%  '
Set-TextValue "E17" '  +0.95%  '
Set-TextValue "D18" '2.482.77'
Set-TextValue "E18" '  +1.65%  '
Set-TextValue "D19" '7.57'
Set-TextValue "E19" '  -1.05%  '
Set-TextValue "D20" '10.97'
Set-TextValue "E20" '  -3.52%  '
Set-TextValue "D21" '350.35'
Set-TextValue "E21" '  -1.44%  '
Set-TextValue "D22" '3.99'
Set-TextValue "E22" '  -0.77%  '
Set-TextValue "E23" '  +0.07%  '
Set-TextValue "D24" '69.03'
Set-TextValue "E25" '  +0.21%  '
Set-TextValue "E26" '  +2.99%  '
Set-TextValue "D27" '9.18'
Set-TextValue "E27" '  +1.95%  '
Set-TextValue "D28" '2.607.68'
Set-TextValue "E28" '  +0.52%  '
Set-TextValue "E29" '  +0.60%  '
Set-TextValue "D30" '0.0₃0912'
Set-TextValue "E30" '  +1.32%  '
Set-TextValue "D31" '507.11'
Set-TextValue "E31" '  -0.70%  '
Set-TextValue "E32" '  -1.19%  '
Set-TextValue "D33" '1.25'
Set-TextValue "E33" '  +1.72%  '
Set-TextValue "E34" '  -0.81%  '
Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  -0.03%  '
Set-TextValue "D36" '162.45'
Set-TextValue "E36" '  +2.44%  '
Set-TextValue "E37" '  -0.60%  '
Set-TextValue "E38" '  +0.64%  '
Set-TextValue "D39" '18.15'
Set-TextValue "E39" '  -1.58%  '
Set-TextValue "E40" '  -0.68%  '
Set-TextValue "E41" '  -0.03%  '
Set-TextValue "E42" '  +1.07%  '
Set-TextValue "E43" '  +1.31%  '
Set-TextValue "D44" '4.84'
Set-TextValue "E44" '  +0.99%  '
Set-TextValue "D45" '2.40'
Set-TextValue "E45" '  +3.30%  '
Set-TextValue "D46" '143.38'
Set-TextValue "E46" '  +1.35%  '
Set-TextValue "D47" '0.0₆0263'
Set-TextValue "E47" '  +3.36%  '
Set-TextValue "E48" '  +0.13%  '
Set-TextValue "D49" '0.516'
Set-TextValue "E49" '  -0.26%  '
Set-TextValue "D50" '0.0738'
Set-TextValue "E50" '  +0.45%  '
Set-TextValue "E51" '  -1.02%  '
